# Metanome_DbSchema.pptx edit: corrected foreign key relation in db schema
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text
#    (09.01.14 -> 12.01.14) on the slide master and every slide layout.
# 2) Reroute the "algorithm" foreign-key connector between the
#    Result/algorithm box and the Execution/algorithm box lower on the
#    slide (5 connector shapes repositioned).

$p = $ppt.ActivePresentation

function Get-PlaceholderShape($shapes, $phType) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $phType) {
            return $sh
        }
    }
    return $null
}

# --- 1) Update the cached date field text everywhere it is cached ---
$newDate = "12.01.14"

$master = $p.SlideMaster
$masterDateShape = Get-PlaceholderShape $master.Shapes 16
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $layoutDateShape = Get-PlaceholderShape $layout.Shapes 16
    if ($layoutDateShape -ne $null) {
        $layoutDateShape.TextFrame.TextRange.Text = $newDate
    }
}

# --- 2) Reposition the connectors that make up the FK relation line ---
$slide = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Gerade Verbindung 73 (id 74): vertical drop from the Result box
$sh74 = Get-ShapeById $slide 74
$sh74.Left = 208.41692913385828
$sh74.Top = 311.9617462834646
$sh74.Width = 0.0
$sh74.Height = 21.64779577559055

# Gerade Verbindung 74 (id 75): short horizontal leg near the target
$sh75 = Get-ShapeById $slide 75
$sh75.Left = 192.39417322834646
$sh75.Top = 333.5046539692913
$sh75.Width = 16.02275590551181
$sh75.Height = 0.0

# Gerade Verbindung 75 (id 76): horizontal leg under the Result box
$sh76 = Get-ShapeById $slide 76
$sh76.Left = 187.34740457480314
$sh76.Top = 311.9617462834646
$sh76.Width = 22.28062992125984
$sh76.Height = 0.0

# Gerade Verbindung 85 (id 86): long vertical spine, extended further down
$sh86 = Get-ShapeById $slide 86
$sh86.Left = 226.8840157480315
$sh86.Top = 97.2048031496063
$sh86.Width = 0.07622047244094488
$sh86.Height = 224.10653543307086

# Gerade Verbindung 92 (id 93): short flipped connector near the spine
$sh93 = Get-ShapeById $slide 93
$sh93.Left = 208.34070866141732
$sh93.Top = 321.31133858267714
$sh93.Width = 18.543307386614174
$sh93.Height = 0.1505511811023622
